$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1224
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 5057
$ws1.Range("F6").Value = 529
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 249
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 693
$ws1.Range("F12").Value = 76

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F3").Value = 20
$ws2.Range("F6").Value = 3

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1224
$ws4.Range("F3").Value = 651
$ws4.Range("F4").Value = 352
$ws4.Range("F7").Value = 5057
$ws4.Range("F8").Value = 529
$ws4.Range("F9").Value = 10
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 249
$ws4.Range("F14").Value = 7
$ws4.Range("F16").Value = 693
$ws4.Range("F18").Value = 76
